$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that no longer exist in the updated dataset
$ws.Rows(6).Delete()
$ws.Rows(6).Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gast"
$ws.Range("C2").Value = "Cckbr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2607106666666667
$ws.Range("H2").Value = 0.782132
$ws.Range("I2").Value = 0.6744482444022483
$ws.Range("J2").Value = 0.6744482444022483
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08226166666666666
$ws.Range("N2").Value = 0.246785
$ws.Range("O2").Value = 0.9036565896241614
$ws.Range("P2").Value = 0.9036565896241615
$ws.Range("Q2").Value = 0.02144649395777778
$ws.Range("R2").Value = 0.19301844562
$ws.Range("S2").Value = 0.6094696004145386
$ws.Range("T2").Value = 0.6094696004145387

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gast"
$ws.Range("C3").Value = "Cckbr"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2607106666666667
$ws.Range("H3").Value = 0.782132
$ws.Range("I3").Value = 0.6744482444022483
$ws.Range("J3").Value = 0.6744482444022483
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.008770333333333333
$ws.Range("N3").Value = 0.026311
$ws.Range("O3").Value = 0.09634341037583853
$ws.Range("P3").Value = 0.09634341037583853
$ws.Range("Q3").Value = 0.002286519450222223
$ws.Range("R3").Value = 0.020578675052
$ws.Range("S3").Value = 0.06497864398770965
$ws.Range("T3").Value = 0.06497864398770965

# Row 4
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Gast"
$ws.Range("C4").Value = "Cckbr"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1258433333333333
$ws.Range("H4").Value = 0.37753
$ws.Range("I4").Value = 0.3255517555977517
$ws.Range("J4").Value = 0.3255517555977517
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.08226166666666666
$ws.Range("N4").Value = 0.246785
$ws.Range("O4").Value = 0.9036565896241614
$ws.Range("P4").Value = 0.9036565896241615
$ws.Range("Q4").Value = 0.01035208233888889
$ws.Range("R4").Value = 0.09316874104999999
$ws.Range("S4").Value = 0.2941869892096228
$ws.Range("T4").Value = 0.2941869892096229

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Gast"
$ws.Range("C5").Value = "Cckbr"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1258433333333333
$ws.Range("H5").Value = 0.37753
$ws.Range("I5").Value = 0.3255517555977517
$ws.Range("J5").Value = 0.3255517555977517
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.008770333333333333
$ws.Range("N5").Value = 0.026311
$ws.Range("O5").Value = 0.09634341037583853
$ws.Range("P5").Value = 0.09634341037583853
$ws.Range("Q5").Value = 0.001103687981111111
$ws.Range("R5").Value = 0.00993319183
$ws.Range("S5").Value = 0.03136476638812889
$ws.Range("T5").Value = 0.03136476638812889
